$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 587; everything from the old row 587
# downward shifts down by one (old 587 -> 588, ..., old 619 -> 620).
$ws.Rows("587:587").Insert()

# Populate the newly inserted row 587 with the new weekly data point.
$ws.Range("A587").Value = 9
$ws.Range("B587").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C587").Value = "Metropolitana"
$ws.Range("D587").Value = 44706
$ws.Range("E587").Value = 13
$ws.Range("F587").Value = 100112040
$ws.Range("G587").Value = "Cilantro"
$ws.Range("H587").Value = "Sin especificar"
$ws.Range("I587").Value = "Primera"
$ws.Range("J587").Value = 175
$ws.Range("K587").Value = 8000
$ws.Range("L587").Value = 9000
$ws.Range("M587").Value = 8429
$ws.Range("N587").Value = '$/docena de atados'
$ws.Range("O587").Value = "Región Metropolitana"
$ws.Range("P587").Value = 2810
$ws.Range("Q587").Value = 3
$ws.Range("R587").Value = "Hortaliza"
